# Daily attendance processing - 2025-12-02 10:56:29
# Normalizes the "Recorded By" (column G) lists on the active sheet so that
# "System" is the leading entry of a two-author list, and otherwise swaps
# the last two recorded authors into their corrected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $count = $parts.Count

    $shouldSwap = $false
    if ($count -gt 2) {
        $shouldSwap = $true
    } elseif ($count -eq 2 -and $parts[0] -ne "System") {
        $shouldSwap = $true
    }

    if ($shouldSwap) {
        $last = $count - 1
        $secondLast = $count - 2
        $tmp = $parts[$last]
        $parts[$last] = $parts[$secondLast]
        $parts[$secondLast] = $tmp

        $cell.Value = [string]::Join(", ", $parts)
    }
}
